$d = $word.ActiveDocument

# The document contains a 5-column practice-problem table ("XX÷Y=" cells).
# Each cell's text is unique within the document, so each value is updated
# via a targeted Find & Replace. Order matters only where a replacement's
# new value happens to equal another cell's old value (20÷2=); that case is
# handled by performing the "20÷2= -> 87÷9=" substitution before the
# "38÷7= -> 20÷2=" substitution, so no double-replacement occurs.

$d.Content.Find.Execute("37÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷5=", 2) | Out-Null
$d.Content.Find.Execute("86÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷4=", 2) | Out-Null
$d.Content.Find.Execute("66÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("65÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=", 2) | Out-Null
$d.Content.Find.Execute("83÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "78÷3=", 2) | Out-Null

$d.Content.Find.Execute("48÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷5=", 2) | Out-Null
$d.Content.Find.Execute("66÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("64÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=", 2) | Out-Null
$d.Content.Find.Execute("22÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2) | Out-Null
$d.Content.Find.Execute("54÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷6=", 2) | Out-Null

$d.Content.Find.Execute("76÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷6=", 2) | Out-Null
$d.Content.Find.Execute("29÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=", 2) | Out-Null
$d.Content.Find.Execute("20÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷9=", 2) | Out-Null
$d.Content.Find.Execute("38÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=", 2) | Out-Null
$d.Content.Find.Execute("21÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "95÷3=", 2) | Out-Null

$d.Content.Find.Execute("90÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=", 2) | Out-Null
$d.Content.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=", 2) | Out-Null
$d.Content.Find.Execute("15÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷2=", 2) | Out-Null
$d.Content.Find.Execute("12÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷6=", 2) | Out-Null
$d.Content.Find.Execute("66÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=", 2) | Out-Null

$d.Content.Find.Execute("42÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=", 2) | Out-Null
$d.Content.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "18÷7=", 2) | Out-Null
$d.Content.Find.Execute("45÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷7=", 2) | Out-Null
$d.Content.Find.Execute("29÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷9=", 2) | Out-Null
$d.Content.Find.Execute("81÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷7=", 2) | Out-Null

Write-Output "Done"
